# Fix typo in the cover-sheet return address: "HMCTS SCSS" -> "HMCTS FPLA"
# (split into two runs, "HMCTS " and "FPLA"), and move the "_GoBack" bookmark
# from the end of the "CM20 9RT" line to the end of the new "HMCTS FPLA" line.

$d = $word.ActiveDocument

# Locate the "HMCTS SCSS" paragraph on the cover sheet.
$hmctsPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "HMCTS SCSS`r") {
        $hmctsPara = $p
        break
    }
}

$start = $hmctsPara.Range.Start

# Replace the paragraph's text, temporarily appending an extra marker
# character so that the bookmark we add below does not land exactly on the
# paragraph-mark boundary (collapsed bookmarks placed there get mis-anchored).
$oldLen = 10   # length of "HMCTS SCSS"
$textRange = $d.Range($start, $start + $oldLen)
$textRange.Text = "HMCTS FPLAX"

# Split "HMCTS " and "FPLA" into two separate runs (same formatting) by
# briefly toggling a character property on the first word and back off.
$firstWordRange = $d.Range($start, $start + 6)
$firstWordRange.Font.Bold = $true
$firstWordRange.Font.Bold = $false

# Move the "_GoBack" bookmark here (collapsed, right after "FPLA").
$newBookmarkPos = $start + 10
$bookmarkRange = $d.Range($newBookmarkPos, $newBookmarkPos)

$oldBookmark = $d.Bookmarks("_GoBack")
$oldBookmark.Delete()

$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# Now remove the temporary "X" marker character.
$markerRange = $d.Range($start + 10, $start + 11)
$markerRange.Text = ""
